# Replace the public-exposure-sites data with a single new entry
# (Melbourne / Nandos) and drop the old extra rows (3-8), leaving only
# the header row and the one data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-stale rows (Black Rock, Hampton, Keysborough, Springvale x2)
$ws.Rows("3:8").Delete()

# Overwrite the remaining data row with the new exposure site details.
# (D2/E2 -- "Case dined at venue" / "new" -- stay as-is.)
$ws.Range("A2").Value = "Melbourne"
$ws.Range("B2").Value = "Nandos  27 Elizabeth Street, Melbourne"
$ws.Range("C2").Value = "01/01/2021 2:00am - 2:30am"

# Match Excel's recalculated "best fit" column widths for the new,
# shorter content (col E is untouched so it keeps its original width).
$ws.Columns("A").ColumnWidth = 8.565104166666666
$ws.Columns("B").ColumnWidth = 31.498697916666668
$ws.Columns("C").ColumnWidth = 23.565104166666668
$ws.Columns("D").ColumnWidth = 15.697916666666666
